$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '30.406.15'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '1.859.22'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = "'235.04"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = "'0.4747"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.10%  '
$ws.Range('D8').Value = "'0.2747"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.10%  '
$ws.Range('D9').Value = "'0.06441"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').Value = '1.886.90'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').Value = "'0.07440"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').Value = "'5.006"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.50%  '
$ws.Range('D14').Value = "'85.81"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').Value = "'0.6356"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.88%  '
$ws.Range('D16').Value = '30.355.03'
$ws.Range('E16').Value = '  -0.62%  '
$ws.Range('D17').Value = "'1.000"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').Value = "'232.27"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.47%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = "'12.81"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.05%  '
$ws.Range('D20').Value = "'0.000007424"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.61%  '
$ws.Range('D21').Value = '2.104.47'
$ws.Range('E21').Value = '  -3.49%  '
$ws.Range('D22').Value = "'1.001"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  -4.45%  '
$ws.Range('D24').Value = "'6.023"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.74%  '
$ws.Range('D25').Value = "'9.238"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('D26').Value = "'166.29"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').Value = "'17.98"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.28%  '
$ws.Range('D28').Value = "'1.898"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').Value = "'0.1035"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.15%  '
$ws.Range('D30').Value = "'1.396"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D31').Value = "'4.162"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.78%  '
$ws.Range('D32').Value = "'3.941"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('D33').Value = "'0.04918"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.50%  '
$ws.Range('D34').Value = "'1.157"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.65%  '
$ws.Range('D35').Value = "'0.7288"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('D36').Value = "'0.9997"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = "'2.696"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').Value = "'0.01897"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.30%  '
$ws.Range('D39').Value = "'2.648"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.51%  '
$ws.Range('D40').Value = "'0.9161"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('D41').Value = "'1.974"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.10%  '
$ws.Range('D42').Value = "'105.98"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').Value = "'0.9999"
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = "'0.4123"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('D45').Value = "'5.591"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range('D46').Value = "'7.139"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.26%  '
$ws.Range('D47').Value = "'61.20"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.51%  '
$ws.Range('D48').Value = "'0.1213"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.46%  '
$ws.Range('D49').Value = "'8.725"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.70%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = "'33.56"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = "'1.410"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.60%  '
